$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting existing rows 51:68 down to 52:69,
# which also extends the used range / dimension to R69.
$ws.Rows(51).Insert()

# Populate the newly inserted row 51 with the new record.
$ws.Range("A51").Value = 11
$ws.Range("B51").Value = "Vega Monumental Concepción"
$ws.Range("C51").Value = "Bíobío"
$ws.Range("D51").Value = 45205
$ws.Range("E51").Value = 8
$ws.Range("F51").Value = 100112026
$ws.Range("G51").Value = "Haba"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 50
$ws.Range("K51").Value = 14000
$ws.Range("L51").Value = 14000
$ws.Range("M51").Value = 14000
$ws.Range("N51").Value = "$/saco 25 kilos"
$ws.Range("O51").Value = "Región Metropolitana"
$ws.Range("P51").Value = 560
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"
